$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Service ID column (A) - BBA9 -> BFA1
$ws.Range("A2").Value = "BFA1"
$ws.Range("A3").Value = "BFA1"
$ws.Range("A4").Value = "BFA1"
$ws.Range("A5").Value = "BFA1"

# Update the Role column (C) - replace iac roles with ia/civil roles
$ws.Range("C2").Value = "caseworker-ia"
$ws.Range("C3").Value = "caseworker-ia-caseofficer"
$ws.Range("C4").Value = "caseworker-civil"
$ws.Range("C5").Value = "caseworker-civil-staff"

# Update the active selection to match the saved view state
$ws.Range("B7").Select()
